$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '67.054.52'
$ws.Range("E2").Value = '  -0.94%  '
$ws.Range("D3").Value = '3.522.94'
$ws.Range("E3").Value = '  +0.51%  '
$ws.Range("E4").Value = '  +0.09%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '609.33'
$ws.Range("E5").Value = '  +0.51%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '147.83'
$ws.Range("E6").Value = '  -2.80%  '
$ws.Range("D7").Value = '3.521.90'
$ws.Range("E7").Value = '  +0.50%  '
$ws.Range("E8").Value = '  +0.03%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.479'
$ws.Range("E9").Value = '  -2.03%  '
$ws.Range("E10").Value = '  -0.66%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '7.81'
$ws.Range("E11").Value = '  +2.42%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.423'
$ws.Range("E12").Value = '  -1.97%  '
$ws.Range("E13").Value = '  -0.24%  '
$ws.Range("D14").Value = '4.116.88'
$ws.Range("E14").Value = '  +0.51%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '31.65'
$ws.Range("E15").Value = '  -2.31%  '
$ws.Range("D16").Value = '3.522.42'
$ws.Range("E16").Value = '  +0.19%  '
$ws.Range("D17").Value = '67.053.81'
$ws.Range("E17").Value = '  -0.77%  '
$ws.Range("E18").Value = '  +0.06%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '10.70'
$ws.Range("E19").Value = '  +8.34%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '6.45'
$ws.Range("E20").Value = '  -1.15%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '15.36'
$ws.Range("E21").Value = '  -1.02%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '437.36'
$ws.Range("E22").Value = '  -2.19%  '
$ws.Range("E23").Value = '  -3.17%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '79.71'
$ws.Range("E24").Value = '  +1.80%  '
$ws.Range("D25").Value = '3.650.89'
$ws.Range("E25").Value = '  +0.20%  '
$ws.Range("E26").Value = '  -0.01%  '
$ws.Range("E27").Value = '  -3.71%  '
$ws.Range("E28").Value = '  -2.83%  '
$ws.Range("E29").Value = '  -4.85%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '2.52'
$ws.Range("E30").Value = '  +0.09%  '
$ws.Range("E31").Value = '  -3.31%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.168'
$ws.Range("E32").Value = '  -1.94%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.999'
$ws.Range("E33").Value = '  -0.66%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '25.43'
$ws.Range("E34").Value = '  -0.80%  '
$ws.Range("D35").Value = '3.517.72'
$ws.Range("E35").Value = '  +0.67%  '
$ws.Range("E36").Value = '  -3.20%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.80'
$ws.Range("E37").Value = '  -3.07%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '8.09'
$ws.Range("E38").Value = '  +1.42%  '
$ws.Range("E39").Value = '  +0.00%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.999'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.0895'
$ws.Range("E41").Value = '  +0.11%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '171.51'
$ws.Range("E42").Value = '  -3.43%  '
$ws.Range("E43").Value = '  -0.04%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.09'
$ws.Range("E44").Value = '  -10.04%  '
$ws.Range("E45").Value = '  +0.45%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.34'
$ws.Range("E46").Value = '  +3.13%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '45.99'
$ws.Range("E47").Value = '  -1.03%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '28.22'
$ws.Range("E48").Value = '  -6.74%  '
$ws.Range("E49").Value = '  -1.81%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.45'
$ws.Range("E50").Value = '  -4.74%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.990'
$ws.Range("E51").Value = '  -0.10%  '
